$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data
$ws.Range("B2").Value = "Sex"
$ws.Range("C2").Value = 48
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 34
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 26.5323113860929
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "Significant"

# A2 carries the same header-row formatting (bold, bordered, centered) as
# the row-1 label cells, so copy that formatting over before setting the
# value (matches s="1" in the target sheet XML).
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0
